$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: create row 78 with the same style as an existing data row (copy row 77 -> row 78)
$ws.Range("A77:F77").Copy($ws.Range("A78:F78"))

# Step 2: write sorted data (LNr ascending) into rows 2..78
$ws.Cells.Item(2, 1).Value = "L_AGEB_1"
$ws.Cells.Item(2, 2).Value = "Q_AGEB"
$ws.Cells.Item(2, 3).Value = "Primärenergieverbrauch"
$ws.Cells.Item(2, 4).Value = "Primary energy consumption (only available in German)"
$ws.Cells.Item(2, 5).Value = "https://ag-energiebilanzen.de/daten-und-fakten/primaerenergieverbrauch/"
$ws.Cells.Item(2, 6).Value = ""

$ws.Cells.Item(3, 1).Value = "L_BB_1"
$ws.Cells.Item(3, 2).Value = "Q_BUNDESBANK"
$ws.Cells.Item(3, 3).Value = "Verschuldung gem. Maastricht-Vertrag - Deutschland - Gesamtstaat"
$ws.Cells.Item(3, 4).Value = "Deficit / surplus and debt level as defined in the Maastricht Treaty/Germany/Debt level/Debt by category/instrument"
$ws.Cells.Item(3, 5).Value = "https://www.bundesbank.de/dynamic/action/de/statistiken/zeitreihen-datenbanken/zeitreihen-datenbank/759778/759778?listId=www_v27_web011_21a"
$ws.Cells.Item(3, 6).Value = "https://www.bundesbank.de/dynamic/action/en/statistics/time-series-databases/time-series-databases/743796/743796?treeAnchor=FINANZEN&statisticType=BBK_ITS"

$ws.Cells.Item(4, 1).Value = "L_BB_2"
$ws.Cells.Item(4, 2).Value = "Q_BUNDESBANK"
$ws.Cells.Item(4, 3).Value = "Studie zur wirtschaftlichen Lage privater Haushalte"
$ws.Cells.Item(4, 4).Value = "Panel on household finances (PHF)"
$ws.Cells.Item(4, 5).Value = "https://www.bundesbank.de/de/bundesbank/forschung/haushaltsstudie"
$ws.Cells.Item(4, 6).Value = "https://www.bundesbank.de/en/bundesbank/research/panel-on-household-finances"

$ws.Cells.Item(5, 1).Value = "L_BFN_1"
$ws.Cells.Item(5, 2).Value = "Q_BFN"
$ws.Cells.Item(5, 3).Value = "Indikator „Artenvielfalt und Landschaftsqualität“"
$ws.Cells.Item(5, 4).Value = "„Biodiversity and landscape quality“ indicator (only available in German)"
$ws.Cells.Item(5, 5).Value = "https://www.bfn.de/indikator-artenvielfalt-und-landschaftsqualitaet"
$ws.Cells.Item(5, 6).Value = ""

$ws.Cells.Item(6, 1).Value = "L_BKA_1"
$ws.Cells.Item(6, 2).Value = "Q_BKA"
$ws.Cells.Item(6, 3).Value = "Polizeiliche Kriminalstatistik"
$ws.Cells.Item(6, 4).Value = "Police Crime Statistics"
$ws.Cells.Item(6, 5).Value = "https://www.bka.de/DE/AktuelleInformationen/StatistikenLagebilder/PolizeilicheKriminalstatistik/pks_node.html"
$ws.Cells.Item(6, 6).Value = "https://www.bka.de/EN/CurrentInformation/PoliceCrimeStatistics/policecrimestatistics_node.html"

$ws.Cells.Item(7, 1).Value = "L_BMDV_1"
$ws.Cells.Item(7, 2).Value = "Q_BMDV"
$ws.Cells.Item(7, 3).Value = "Breitbandatlas"
$ws.Cells.Item(7, 4).Value = " Broadband Atlas (only available in German)"
$ws.Cells.Item(7, 5).Value = "https://www.bundesnetzagentur.de/DE/Fachthemen/Telekommunikation/Breitband/breitbandatlas/start.html"
$ws.Cells.Item(7, 6).Value = ""

$ws.Cells.Item(8, 1).Value = "L_BMEL_1"
$ws.Cells.Item(8, 2).Value = "Q_BMEL"
$ws.Cells.Item(8, 3).Value = "Nährstoffbilanzen und Düngemittel (Tabelle MBT-0111260-0000)"
$ws.Cells.Item(8, 4).Value = "Nutrient balances and fertilizers (table MBT-0111260-0000; only available in German)"
$ws.Cells.Item(8, 5).Value = "https://www.bmel-statistik.de/landwirtschaft/statistischer-monatsbericht-des-bmel-kapitel-a-landwirtschaft/"
$ws.Cells.Item(8, 6).Value = "https://www.bmel-statistik.de/landwirtschaft/statistischer-monatsbericht-des-bmel-kapitel-a-landwirtschaft/"

$ws.Cells.Item(9, 1).Value = "L_BMEL_2"
$ws.Cells.Item(9, 2).Value = "Q_BMEL"
$ws.Cells.Item(9, 3).Value = "Betriebe und Flächen des ökologischen Landbaus in Deutschland"
$ws.Cells.Item(9, 4).Value = "Agricultural holdings total and holdings with organic farming"
$ws.Cells.Item(9, 5).Value = "https://www.bmel.de/DE/themen/landwirtschaft/oekologischer-landbau/oekologischer-landbau_node"
$ws.Cells.Item(9, 6).Value = "https://www.bmel.de/EN/topics/farming/organic-farming/organic-farming_node.html"

$ws.Cells.Item(10, 1).Value = "L_BMEL_3"
$ws.Cells.Item(10, 2).Value = "Q_BMEL"
$ws.Cells.Item(10, 3).Value = "Ökologischer Landbau in Deutschland"
$ws.Cells.Item(10, 4).Value = "Organic farming in Germany (only available in German)"
$ws.Cells.Item(10, 5).Value = "https://www.bmel.de/SharedDocs/Downloads/DE/Broschueren/OekolandbauDeutschland.pdf?__blob=publicationFile&v=14"
$ws.Cells.Item(10, 6).Value = ""

$ws.Cells.Item(11, 1).Value = "L_BMZ_1"
$ws.Cells.Item(11, 2).Value = "Q_BMZ"
$ws.Cells.Item(11, 3).Value = "Zahlen und Fakten der deutschen Entwicklungszusammenarbeit"
$ws.Cells.Item(11, 4).Value = "Facts and figures of German development cooperation"
$ws.Cells.Item(11, 5).Value = "https://www.bmz.de/de/ministerium/zahlen-fakten"
$ws.Cells.Item(11, 6).Value = "https://www.bmz.de/en/ministry/facts-figures"

$ws.Cells.Item(12, 1).Value = "L_BMZ_2"
$ws.Cells.Item(12, 2).Value = "Q_BMZ"
$ws.Cells.Item(12, 3).Value = "Elmau Progress Report 2022 (nur auf Englisch verfügbar)"
$ws.Cells.Item(12, 4).Value = "Elmau Progress Report 2022"
$ws.Cells.Item(12, 5).Value = "https://www.bmz.de/de/aktuelles/g7-praesidentschaft/elmau-progress-report-2022"
$ws.Cells.Item(12, 6).Value = "https://www.bmz.de/resource/blob/116286/4e535985a12d0701e63b25ca9f4d0fb8/2022-07-01-elmau-progress-report-2022-data.pdf"

$ws.Cells.Item(13, 1).Value = "L_BZGA_1"
$ws.Cells.Item(13, 2).Value = "Q_BZGA"
$ws.Cells.Item(13, 3).Value = "Suchtprävention"
$ws.Cells.Item(13, 4).Value = "Addiction Prevention (only available in German)"
$ws.Cells.Item(13, 5).Value = "https://www.bzga.de/presse/daten-und-fakten/suchtpraevention/"
$ws.Cells.Item(13, 6).Value = "https://www.bzga.de/presse/daten-und-fakten/suchtpraevention/"

$ws.Cells.Item(14, 1).Value = "L_BZM_1"
$ws.Cells.Item(14, 2).Value = "Q_BMZ"
$ws.Cells.Item(14, 3).Value = "Klimafinanzierung"
$ws.Cells.Item(14, 4).Value = "Climate financing"
$ws.Cells.Item(14, 5).Value = "https://www.bmz.de/de/entwicklungspolitik/klimawandel-und-entwicklung/klimafinanzierung"
$ws.Cells.Item(14, 6).Value = ""

$ws.Cells.Item(15, 1).Value = "L_DSTTS_10"
$ws.Cells.Item(15, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(15, 3).Value = "Internationale Bildungsindikatoren im Ländervergleich"
$ws.Cells.Item(15, 4).Value = "International education indicators in country comparison (only available in German)"
$ws.Cells.Item(15, 5).Value = "https://www.destatis.de/DE/Themen/Gesellschaft-Umwelt/Bildung-Forschung-Kultur/Bildungsstand/_inhalt.html#sprg233662"
$ws.Cells.Item(15, 6).Value = "https://www.destatis.de/DE/Themen/Gesellschaft-Umwelt/Bildung-Forschung-Kultur/Bildungsstand/_inhalt.html#sprg233662"

$ws.Cells.Item(16, 1).Value = "L_DSTTS_11"
$ws.Cells.Item(16, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(16, 3).Value = "Kinder und tätige Personen in Tageseinrichtungen und in öffentlich geförderter Kindertagespflege"
$ws.Cells.Item(16, 4).Value = "Children and persons working in day care facilities and in publicly funded day care for children (only available in German)"
$ws.Cells.Item(16, 5).Value = "https://www.destatis.de/DE/Themen/Gesellschaft-Umwelt/Soziales/Kindertagesbetreuung/_inhalt.html#sprg234640"
$ws.Cells.Item(16, 6).Value = ""

$ws.Cells.Item(17, 1).Value = "L_DSTTS_12"
$ws.Cells.Item(17, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(17, 3).Value = "VGR des Bundes - Bruttoanlageinvestitionen (nominal/preisbereinigt) – GENESIS online 81000-0023"
$ws.Cells.Item(17, 4).Value = "National accounts - Gross fixed capital formation (nominal/price-adjusted) – GENESIS online 81000-0023"
$ws.Cells.Item(17, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=81000-0023&bypass=true&levelindex=0&levelid=1660823284613&language=de"
$ws.Cells.Item(17, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=81000-0023&bypass=true&levelindex=0&levelid=1660823284613&language=en"

$ws.Cells.Item(18, 1).Value = "L_DSTTS_13"
$ws.Cells.Item(18, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(18, 3).Value = "VGR des Bundes - Bruttowertschöpfung, Bruttoinlandsprodukt (nominal/preisbereinigt) – GENESIS online 81000-0001"
$ws.Cells.Item(18, 4).Value = "National accounts - Gross value added, gross domestic product (nominal/price-adjusted) – GENESIS online 81000-0001"
$ws.Cells.Item(18, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code= 81000-0001 &bypass=true&levelindex=0&levelid=1660823284613&language=de"
$ws.Cells.Item(18, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code= 81000-0001 &bypass=true&levelindex=0&levelid=1660823284613&language=en"

$ws.Cells.Item(19, 1).Value = "L_DSTTS_14"
$ws.Cells.Item(19, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(19, 3).Value = "Durchschnittliche Bevölkerung – GENESIS online 12411-0040"
$ws.Cells.Item(19, 4).Value = "Average population – GENESIS online 12411-0040"
$ws.Cells.Item(19, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code= 12411-0040 &bypass=true&levelindex=0&levelid=1660823284613&language=de"
$ws.Cells.Item(19, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code= 12411-0040 &bypass=true&levelindex=0&levelid=1660823284613&language=en"

$ws.Cells.Item(20, 1).Value = "L_DSTTS_15"
$ws.Cells.Item(20, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(20, 3).Value = "Rückgerechnete und fortgeschriebene Bevölkerung auf Grundlage des Zensus 2011"
$ws.Cells.Item(20, 4).Value = "Back-calculated and updated population based on the 2011 census"
$ws.Cells.Item(20, 5).Value = "https://www.destatis.de/DE/Themen/Gesellschaft-Umwelt/Bevoelkerung/Bevoelkerungsstand/_inhalt.html#sprg233540"
$ws.Cells.Item(20, 6).Value = ""

$ws.Cells.Item(21, 1).Value = "L_DSTTS_16"
$ws.Cells.Item(21, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(21, 3).Value = "Absolventen und Abgänger: Deutschland – GENESIS online 21111-0004"
$ws.Cells.Item(21, 4).Value = "School leavers (graduates and dropouts): Germany – GENESIS online 21111-0004"
$ws.Cells.Item(21, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21111-0004&bypass=true&levelindex=1&levelid=1660810680251&language=de"
$ws.Cells.Item(21, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21111-0004&bypass=true&levelindex=0&levelid=1660823284613&language=en"

$ws.Cells.Item(22, 1).Value = "L_DSTTS_17"
$ws.Cells.Item(22, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(22, 3).Value = "Entwicklungszusammenarbeit"
$ws.Cells.Item(22, 4).Value = "Development cooperation"
$ws.Cells.Item(22, 5).Value = "https://www.destatis.de/DE/Themen/Staat/Oeffentliche-Finanzen/Entwicklungszusammenarbeit/_inhalt.html"
$ws.Cells.Item(22, 6).Value = "https://www.destatis.de/EN/Themes/Government/Public-Finance/Development-Cooperation/_node.html"

$ws.Cells.Item(23, 1).Value = "L_DSTTS_18"
$ws.Cells.Item(23, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(23, 3).Value = "Aus- und Einfuhr (Außenhandel) – GENESIS online 51000-0001"
$ws.Cells.Item(23, 4).Value = "Exports and imports (foreign trade) – GENESIS online 51000-0001"
$ws.Cells.Item(23, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=51000-0007&bypass=true&levelindex=1&levelid=1669021022626&language=de"
$ws.Cells.Item(23, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=51000-0007&bypass=true&levelindex=1&levelid=1669021022626&language=en"

$ws.Cells.Item(24, 1).Value = "L_DSTTS_19"
$ws.Cells.Item(24, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(24, 3).Value = "Unbereinigter Gender Pay Gap nach Bundesländern"
$ws.Cells.Item(24, 4).Value = "Unadjusted Gender Pay Gap by Länder"
$ws.Cells.Item(24, 5).Value = "https://www.destatis.de/DE/Themen/Arbeit/Verdienste/Verdienste-Verdienstunterschiede/Tabellen/ugpg-02-bundeslaender-ab-2014.html"
$ws.Cells.Item(24, 6).Value = "https://www.destatis.de/EN/Themes/Labour/Earnings/Earnings-Earnings-Differences/Tables/ugpg-02-by-laender-at2014.html"

$ws.Cells.Item(25, 1).Value = "L_DSTTS_2"
$ws.Cells.Item(25, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(25, 3).Value = "Bruttoinlandsprodukt – GENESIS online 81000-0001"
$ws.Cells.Item(25, 4).Value = "National accounts - Gross value added, gross domestic product – GENESIS online 81000-0001"
$ws.Cells.Item(25, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=81000-0001&bypass=true&levelindex=0&levelid=1660822010108&language=de"
$ws.Cells.Item(25, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=81000-0001&bypass=true&levelindex=1&levelid=1660802268437&language=en"

$ws.Cells.Item(26, 1).Value = "L_DSTTS_20"
$ws.Cells.Item(26, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(26, 3).Value = "Väterbeteiligung – GENESIS online 22922-0011"
$ws.Cells.Item(26, 4).Value = "Participation rate of fathers"
$ws.Cells.Item(26, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=22922-0011&bypass=true&levelindex=0&levelid=1660642440197#abreadcrumb&language=de"
$ws.Cells.Item(26, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=22922-0011&bypass=true&levelindex=0&levelid=1660813986805#abreadcrumb&language=en"

$ws.Cells.Item(27, 1).Value = "L_DSTTS_21"
$ws.Cells.Item(27, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(27, 3).Value = "Absolventen und Abgänger: Bundesländer – GENESIS online 21111-0013"
$ws.Cells.Item(27, 4).Value = "School leavers (graduates and dropouts): Länder – GENESIS online 21111-0013"
$ws.Cells.Item(27, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21111-0013&bypass=true&levelindex=0&levelid=1660823504838&language=de"
$ws.Cells.Item(27, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21111-0013&bypass=true&levelindex=1&levelid=1660810680251&language=en"

$ws.Cells.Item(28, 1).Value = "L_DSTTS_22"
$ws.Cells.Item(28, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(28, 3).Value = "Anstieg der Siedlungs- und Verkehrsfläche"
$ws.Cells.Item(28, 4).Value = "Increase in settlement and transport area (only available in German)"
$ws.Cells.Item(28, 5).Value = "https://www.destatis.de/DE/Themen/Branchen-Unternehmen/Landwirtschaft-Forstwirtschaft-Fischerei/Flaechennutzung/Tabellen/anstieg-suv2.html"
$ws.Cells.Item(28, 6).Value = ""

$ws.Cells.Item(29, 1).Value = "L_DSTTS_23"
$ws.Cells.Item(29, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(29, 3).Value = "Karten zur Statistik der Kindertagesbetreuung"
$ws.Cells.Item(29, 4).Value = "Child day care statistics maps (only available in German)"
$ws.Cells.Item(29, 5).Value = "https://www.destatis.de/DE/Themen/Gesellschaft-Umwelt/Soziales/Kindertagesbetreuung/kindertagesbetreuung-karte.html;#karte3"
$ws.Cells.Item(29, 6).Value = ""

$ws.Cells.Item(30, 1).Value = "L_DSTTS_24"
$ws.Cells.Item(30, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(30, 3).Value = "Überbelastung durch Wohnkosten"
$ws.Cells.Item(30, 4).Value = "Housing cost overburden"
$ws.Cells.Item(30, 5).Value = "https://www.destatis.de/Europa/DE/Thema/Bevoelkerung-Arbeit-Soziales/Soziales-Lebensbedingungen/Wohnkosten.html"
$ws.Cells.Item(30, 6).Value = "https://www.destatis.de/Europa/EN/Topic/Population-Labour-Social-Issues/Social-issues-living-conditions/_node.html;jsessionid=B340DD00C6EEDC7477B2AD2B54E4BC40.live731#587120"

$ws.Cells.Item(31, 1).Value = "L_DSTTS_25"
$ws.Cells.Item(31, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(31, 3).Value = "Umweltökonomische Gesamtrechnungen"
$ws.Cells.Item(31, 4).Value = "Environmental Economic Accounting"
$ws.Cells.Item(31, 5).Value = "https://www.destatis.de/DE/Themen/Gesellschaft-Umwelt/Umwelt/UGR/_inhalt.html"
$ws.Cells.Item(31, 6).Value = "https://www.destatis.de/EN/Themes/Society-Environment/Environment/Environmental-Economic-Accounting/_node.html"

$ws.Cells.Item(32, 1).Value = "L_DSTTS_27"
$ws.Cells.Item(32, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(32, 3).Value = "Gesamtrohstoffproduktivität und ihre Komponenten"
$ws.Cells.Item(32, 4).Value = "Raw material input productivity and ist components"
$ws.Cells.Item(32, 5).Value = "https://www.destatis.de/DE/Themen/Gesellschaft-Umwelt/Umwelt/UGR/rohstoffe-materialfluesse-wasser/Tabellen/gesamtrohstoff-produktivitaet.html"
$ws.Cells.Item(32, 6).Value = "https://www.destatis.de/EN/Themes/Society-Environment/Environment/Environmental-Economic-Accounting/raw-material-flows-water/Tables/total-raw-material-productivity.html"

$ws.Cells.Item(33, 1).Value = "L_DSTTS_3"
$ws.Cells.Item(33, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(33, 3).Value = "VGR des Bundes - Einnahmen und Ausgaben sowie – GENESIS online 81000-0031`nFinanzierungssaldo des Staates – GENESIS online 81000-0031"
$ws.Cells.Item(33, 4).Value = "National accounts - Revenue, expenditure, net lending/net – GENESIS online 81000-0031`nborrowing of general government – GENESIS online 81000-0031"
$ws.Cells.Item(33, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=81000-0031&bypass=true&levelindex=1&levelid=1660802268437&language=de"
$ws.Cells.Item(33, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=81000-0031&bypass=true&levelindex=1&levelid=1660802268437&language=en"

$ws.Cells.Item(34, 1).Value = "L_DSTTS_4"
$ws.Cells.Item(34, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(34, 3).Value = "Interne Ausgaben und Personal für Forschung und Entwicklung (Bund) – GENESIS online 21821-0001"
$ws.Cells.Item(34, 4).Value = "Research and development expenditure and staff (Germany) – GENESIS online 21821-0001"
$ws.Cells.Item(34, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21821-0001&bypass=true&levelindex=1&levelid=1622107294362&language=de"
$ws.Cells.Item(34, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21821-0001&bypass=true&levelindex=0&levelid=1660726117256&language=en"

$ws.Cells.Item(35, 1).Value = "L_DSTTS_5"
$ws.Cells.Item(35, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(35, 3).Value = "Einkommensverteilung (Nettoäquivalenzeinkommen) in Deutschland"
$ws.Cells.Item(35, 4).Value = "Income distribution (equivalised net income)"
$ws.Cells.Item(35, 5).Value = "https://www.destatis.de/DE/Themen/Gesellschaft-Umwelt/Einkommen-Konsum-Lebensbedingungen/Lebensbedingungen-Armutsgefaehrdung/Tabellen/einkommensverteilung-mz-silc.html"
$ws.Cells.Item(35, 6).Value = "https://www.destatis.de/EN/Themes/Society-Environment/Income-Consumption-Living-Conditions/Living-Conditions-Risk-Poverty/Tables/income-distribution-mz-silc.html"

$ws.Cells.Item(36, 1).Value = "L_DSTTS_7"
$ws.Cells.Item(36, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(36, 3).Value = "Agrarstrukturerhebung - Betriebe mit ökologischem Landbau"
$ws.Cells.Item(36, 4).Value = "Agrarstrukturerhebung - farms with organic agriculture (only available in German)"
$ws.Cells.Item(36, 5).Value = "https://www.destatis.de/DE/Themen/Branchen-Unternehmen/Landwirtschaft-Forstwirtschaft-Fischerei/Landwirtschaftliche-Betriebe/_inhalt.html#sprg239572"
$ws.Cells.Item(36, 6).Value = "https://www.destatis.de/DE/Themen/Branchen-Unternehmen/Landwirtschaft-Forstwirtschaft-Fischerei/Landwirtschaftliche-Betriebe/_inhalt.html#sprg239572"

$ws.Cells.Item(37, 1).Value = "L_DSTTS_8"
$ws.Cells.Item(37, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(37, 3).Value = "Unbereinigter Gender Pay Gap"
$ws.Cells.Item(37, 4).Value = "Unadjusted gender pay gap"
$ws.Cells.Item(37, 5).Value = "https://www.destatis.de/DE/Themen/Arbeit/Verdienste/Verdienste-Verdienstunterschiede/Tabellen/ugpg-01-gebietsstand.html"
$ws.Cells.Item(37, 6).Value = "https://www.destatis.de/EN/Themes/Labour/Earnings/Earnings-Earnings-Differences/Tables/ugpg-01-by-territory-gpg.html"

$ws.Cells.Item(38, 1).Value = "L_DSTTS_9"
$ws.Cells.Item(38, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(38, 3).Value = "Interne Ausgaben und Personal für Forschung und Entwicklung (Länder) – GENESIS online 21821-0002"
$ws.Cells.Item(38, 4).Value = "Research and development expenditure and staff (Länder) – GENESIS online 21821-0002"
$ws.Cells.Item(38, 5).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21821-0002&bypass=true&levelindex=0&levelid=1660726117256&language=de"
$ws.Cells.Item(38, 6).Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21821-0002&bypass=true&levelindex=1&levelid=1623135114747&language=en"

$ws.Cells.Item(39, 1).Value = "L_EE_1"
$ws.Cells.Item(39, 2).Value = "Q_AGEESTAT"
$ws.Cells.Item(39, 3).Value = "Entwicklung der erneuerbaren Energien in Deutschland"
$ws.Cells.Item(39, 4).Value = "Development of renewable energies in Germany"
$ws.Cells.Item(39, 5).Value = "https://www.erneuerbare-energien.de/EE/Navigation/DE/Service/Erneuerbare_Energien_in_Zahlen/Zeitreihen/zeitreihen.html"
$ws.Cells.Item(39, 6).Value = "https://www.erneuerbare-energien.de/EE/Navigation/DE/Service/Erneuerbare_Energien_in_Zahlen/Zeitreihen/zeitreihen.html"

$ws.Cells.Item(40, 1).Value = "L_EMAS_1"
$ws.Cells.Item(40, 2).Value = "Q_DIHK"
$ws.Cells.Item(40, 3).Value = "EMAS-Statistiken"
$ws.Cells.Item(40, 4).Value = "EMAS statistics (only available in German)"
$ws.Cells.Item(40, 5).Value = "https://www.emas.de/statistiken/"
$ws.Cells.Item(40, 6).Value = "https://www.emas.de/statistiken/"

$ws.Cells.Item(41, 1).Value = "L_ERSTT_1"
$ws.Cells.Item(41, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(41, 3).Value = "Eurostat Datenbank"
$ws.Cells.Item(41, 4).Value = "Eurostat database"
$ws.Cells.Item(41, 5).Value = "https://ec.europa.eu/eurostat/de/data/database"
$ws.Cells.Item(41, 6).Value = "https://ec.europa.eu/eurostat/web/main/data/database"

$ws.Cells.Item(42, 1).Value = "L_ERSTT_10"
$ws.Cells.Item(42, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(42, 3).Value = "Rate der erheblichen materiellen und sozialen Deprivation - Eurostat-Tabelle [ilc_mdsd11]"
$ws.Cells.Item(42, 4).Value = "Severe material and social deprivation rate - Eurostat table [ilc_mdsd11]"
$ws.Cells.Item(42, 5).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDSD11__custom_3696252/default/table?lang=de"
$ws.Cells.Item(42, 6).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDSD11__custom_3696252/default/table?lang=en"

$ws.Cells.Item(43, 1).Value = "L_ERSTT_11"
$ws.Cells.Item(43, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(43, 3).Value = "Quote der Überbelastung durch Wohnkosten - Eurostat-Tabelle [ilc_lvho07a ]"
$ws.Cells.Item(43, 4).Value = "Housing cost overburden rate - Eurostat table [ilc_lvho07a ]"
$ws.Cells.Item(43, 5).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_LVHO07A/default/table?lang=de&category=livcon.ilc.ilc_lv.ilc_lvho.ilc_lvho_hc"
$ws.Cells.Item(43, 6).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_LVHO07A/default/table?category=livcon.ilc.ilc_lv.ilc_lvho.ilc_lvho_hc"

$ws.Cells.Item(44, 1).Value = "L_ERSTT_12"
$ws.Cells.Item(44, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(44, 3).Value = "Für ökologische Landwirtschaft genutzte Fläche"
$ws.Cells.Item(44, 4).Value = "Area under organic farming"
$ws.Cells.Item(44, 5).Value = "https://ec.europa.eu/eurostat/databrowser/view/sdg_02_40/default/table?lang=de"
$ws.Cells.Item(44, 6).Value = "https://ec.europa.eu/eurostat/databrowser/view/sdg_02_40/default/table?lang=en"

$ws.Cells.Item(45, 1).Value = "L_ERSTT_2"
$ws.Cells.Item(45, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(45, 3).Value = "Erwerbstätigenquoten nach Geschlecht, Alter und Staatsangehörigkeit"
$ws.Cells.Item(45, 4).Value = "Employment rates by sex, age and citizenship"
$ws.Cells.Item(45, 5).Value = "https://appsso.eurostat.ec.europa.eu/nui/show.do?dataset=lfsa_ergan&lang=de"
$ws.Cells.Item(45, 6).Value = "https://appsso.eurostat.ec.europa.eu/nui/show.do?dataset=lfsa_ergan&lang=en"

$ws.Cells.Item(46, 1).Value = "L_ERSTT_3"
$ws.Cells.Item(46, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(46, 3).Value = "Gini-Koeffizient des verfügbaren Äquivalenzeinkommens vor Sozialleistungen - EU-SILC Erhebung - Eurostat-Tabelle  [ilc_di12c ]"
$ws.Cells.Item(46, 4).Value = "Gini coefficient of equivalised disposable income before social transfers - Eurostat table  [ilc_di12c ]"
$ws.Cells.Item(46, 5).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12C/default/table?lang=de&category=livcon.ilc.ilc_ie.ilc_iei"
$ws.Cells.Item(46, 6).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12C/default/table?lang=en&category=livcon.ilc.ilc_ie.ilc_iei"

$ws.Cells.Item(47, 1).Value = "L_ERSTT_5"
$ws.Cells.Item(47, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(47, 3).Value = "STECF - Berichte (nur auf Englisch verfügbar)"
$ws.Cells.Item(47, 4).Value = "STECF - Reports"
$ws.Cells.Item(47, 5).Value = ""
$ws.Cells.Item(47, 6).Value = "https://stecf.jrc.ec.europa.eu/reports/cfp-monitoring"

$ws.Cells.Item(48, 1).Value = "L_ERSTT_6"
$ws.Cells.Item(48, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(48, 3).Value = "Gini-Koeffizient des verfügbaren Äquivalenzeinkommens  - EU-SILC Erhebung - Eurostat-Tabelle  [ilc_di12 ]"
$ws.Cells.Item(48, 4).Value = "Gini coefficient of equivalised disposable income - Eurostat table  [ilc_di12c]"
$ws.Cells.Item(48, 5).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12/default/table?lang=de&category=livcon.ilc.ilc_ie.ilc_iei"
$ws.Cells.Item(48, 6).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12/default/table?lang=en&category=livcon.ilc.ilc_ie.ilc_iei"

$ws.Cells.Item(49, 1).Value = "L_ERSTT_7"
$ws.Cells.Item(49, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(49, 3).Value = "Rate der materiellen und sozialen Deprivation - Eurostat-Tabelle [ilc_mdsd07]"
$ws.Cells.Item(49, 4).Value = "Material and social deprivation rate - Eurostat table [ilc_mdsd07]"
$ws.Cells.Item(49, 5).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDSD07/default/table?lang=de"
$ws.Cells.Item(49, 6).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDSD07/default/table?lang=en"

$ws.Cells.Item(50, 1).Value = "L_ERSTT_8"
$ws.Cells.Item(50, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(50, 3).Value = "Rate der erheblichen materiellen Deprivation - Eurostat-Tabelle [ilc_mddd11]"
$ws.Cells.Item(50, 4).Value = "Severe material deprivation rate - Eurostat table [ilc_mddd11]"
$ws.Cells.Item(50, 5).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDDD11/default/table?lang=de&category=livcon.ilc.ilc_md.ilc_mddd"
$ws.Cells.Item(50, 6).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDDD11/default/table?lang=en&category=livcon.ilc.ilc_md.ilc_mddd"

$ws.Cells.Item(51, 1).Value = "L_ERSTT_9"
$ws.Cells.Item(51, 2).Value = "Q_EUROSTAT"
$ws.Cells.Item(51, 3).Value = "Bevölkerung nach Bildungsabschluss (nur Tertiärbereich) - Eurostat-Tabelle [edat_lfs_9912]"
$ws.Cells.Item(51, 4).Value = "Population by educational attainment level (tertiary education only) - Eurostat table [edat_lfs_9912]"
$ws.Cells.Item(51, 5).Value = "https://ec.europa.eu/eurostat/databrowser/view/EDAT_LFS_9912__custom_3201014/default/table?lang=de"
$ws.Cells.Item(51, 6).Value = "https://ec.europa.eu/eurostat/databrowser/view/EDAT_LFS_9912__custom_3201014/default/table?lang=en"

$ws.Cells.Item(52, 1).Value = "L_EZB_1"
$ws.Cells.Item(52, 2).Value = "Q_EZB"
$ws.Cells.Item(52, 3).Value = "Household Finance and Consumption Network (HFCN)"
$ws.Cells.Item(52, 4).Value = "Household Finance and Consumption Network (HFCN)"
$ws.Cells.Item(52, 5).Value = "https://www.ecb.europa.eu/pub/economic-research/research-networks/html/researcher_hfcn.en.html"
$ws.Cells.Item(52, 6).Value = "https://www.ecb.europa.eu/pub/economic-research/research-networks/html/researcher_hfcn.en.html"

$ws.Cells.Item(53, 1).Value = "L_FIDAR_1"
$ws.Cells.Item(53, 2).Value = "Q_FIDAR"
$ws.Cells.Item(53, 3).Value = "Studie zum WoB-Index"
$ws.Cells.Item(53, 4).Value = "Study on the WoB Index (only available in German)"
$ws.Cells.Item(53, 5).Value = "https://www.fidar.de/wob-indizes-studien/wob-index-185/studie-zum-wob-index-185.html"
$ws.Cells.Item(53, 6).Value = "https://www.fidar.de/wob-indizes-studien/wob-index-185/studie-zum-wob-index-185.html"

$ws.Cells.Item(54, 1).Value = "L_GBE_1"
$ws.Cells.Item(54, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(54, 3).Value = "Vorzeitige Sterblichkeit"
$ws.Cells.Item(54, 4).Value = "Premature mortality"
$ws.Cells.Item(54, 5).Value = "https://www.gbe-bund.de/gbe/pkg_isgbe5.prc_menu_olap?p_uid=gast&p_aid=52889592&p_sprache=D&p_help=2&p_indnr=562&p_indsp=3194&p_ityp=H&p_fid="
$ws.Cells.Item(54, 6).Value = "https://www.gbe-bund.de/gbe/pkg_isgbe5.prc_menu_olap?p_uid=gast&p_aid=80721089&p_sprache=E&p_help=2&p_indnr=562&p_version=1&p_ansnr=95662760"

$ws.Cells.Item(55, 1).Value = "L_GIZ_1"
$ws.Cells.Item(55, 2).Value = "Q_BMZ"
$ws.Cells.Item(55, 3).Value = "Mitglieder des Bündnises für nachhaltige Textilien"
$ws.Cells.Item(55, 4).Value = "Members of the Partnership for Sustainable Textiles"
$ws.Cells.Item(55, 5).Value = "https://www.textilbuendnis.com/uebersicht/"
$ws.Cells.Item(55, 6).Value = "https://www.textilbuendnis.com/en/uebersicht/"

$ws.Cells.Item(56, 1).Value = "L_IFEU_1"
$ws.Cells.Item(56, 2).Value = "Q_IFEU"
$ws.Cells.Item(56, 3).Value = "TREMOD"
$ws.Cells.Item(56, 4).Value = "TREMOD"
$ws.Cells.Item(56, 5).Value = "https://www.ifeu.de/methoden-tools/modelle/tremod/"
$ws.Cells.Item(56, 6).Value = "https://www.ifeu.de/en/methods-tools/models/tremod/"

$ws.Cells.Item(57, 1).Value = "L_LAK_1"
$ws.Cells.Item(57, 2).Value = "Q_LAKEB"
$ws.Cells.Item(57, 3).Value = "Energieindikatoren allgemein"
$ws.Cells.Item(57, 4).Value = "Energy indicators general (only available in German)"
$ws.Cells.Item(57, 5).Value = "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i100"
$ws.Cells.Item(57, 6).Value = "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i100"

$ws.Cells.Item(58, 1).Value = "L_LAK_2"
$ws.Cells.Item(58, 2).Value = "Q_LAKEB"
$ws.Cells.Item(58, 3).Value = "Energieindikatoren: Erneuerbare Energieträger"
$ws.Cells.Item(58, 4).Value = "Energy indicators: Renewable energy sources (only available in German)"
$ws.Cells.Item(58, 5).Value = "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i290"
$ws.Cells.Item(58, 6).Value = "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i290"

$ws.Cells.Item(59, 1).Value = "L_LAK_3"
$ws.Cells.Item(59, 2).Value = "Q_LAKEB"
$ws.Cells.Item(59, 3).Value = "Anteil Strom aus erneuerbaren Energiequellen am Bruttostromverbrauch"
$ws.Cells.Item(59, 4).Value = "Share of electricity from renewable energy sources in gross electricity consumption (only available in German)"
$ws.Cells.Item(59, 5).Value = "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i200"
$ws.Cells.Item(59, 6).Value = "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i200"

$ws.Cells.Item(60, 1).Value = "L_LANUV_2"
$ws.Cells.Item(60, 2).Value = "Q_LIKI"
$ws.Cells.Item(60, 3).Value = "Flächenverbrauch"
$ws.Cells.Item(60, 4).Value = "Land consumption (only available in German)"
$ws.Cells.Item(60, 5).Value = "https://www.liki.nrw.de/ressourcen-und-effizienz/d1-flaechenverbrauch"
$ws.Cells.Item(60, 6).Value = ""

$ws.Cells.Item(61, 1).Value = "L_LIKI_1"
$ws.Cells.Item(61, 2).Value = "Q_LIKI"
$ws.Cells.Item(61, 3).Value = "Nitrat im Grundwasser"
$ws.Cells.Item(61, 4).Value = "Nitrate in groundwater (only available in German)"
$ws.Cells.Item(61, 5).Value = "https://www.liki.nrw.de/umwelt-und-gesundheit/c5-nitrat-im-grundwasser"
$ws.Cells.Item(61, 6).Value = ""

$ws.Cells.Item(62, 1).Value = "L_OECD_1"
$ws.Cells.Item(62, 2).Value = "Q_OECD"
$ws.Cells.Item(62, 3).Value = "OECD - Gesamtströme nach Geber (nicht auf Deutsch verfügbar)"
$ws.Cells.Item(62, 4).Value = "OECD - Total flows by donor"
$ws.Cells.Item(62, 5).Value = ""
$ws.Cells.Item(62, 6).Value = "https://stats.oecd.org/Index.aspx?DataSetCode=TABLE1"

$ws.Cells.Item(63, 1).Value = "L_RDB_1"
$ws.Cells.Item(63, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(63, 3).Value = "Bevölkerung"
$ws.Cells.Item(63, 4).Value = "Population (only available in German)"
$ws.Cells.Item(63, 5).Value = "https://www.regionalstatistik.de/genesis/online?operation=previous&levelindex=0&step=0&titel=Tabellenaufbau&levelid=1668672879939&acceptscookies=false#abreadcrumb"
$ws.Cells.Item(63, 6).Value = ""

$ws.Cells.Item(64, 1).Value = "L_RDB_2"
$ws.Cells.Item(64, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(64, 3).Value = "Siedlungsfläche"
$ws.Cells.Item(64, 4).Value = "Settlement area (only available in German)"
$ws.Cells.Item(64, 5).Value = "https://www.regionalstatistik.de/genesis//online?operation=table&code=33111-02-01-4&bypass=true&levelindex=1&levelid=1668501420853#abreadcrumb"
$ws.Cells.Item(64, 6).Value = ""

$ws.Cells.Item(65, 1).Value = "L_RDB_3"
$ws.Cells.Item(65, 2).Value = "Q_DESTATIS"
$ws.Cells.Item(65, 3).Value = "Verkehrsfläche"
$ws.Cells.Item(65, 4).Value = "Transport area (only available in German)"
$ws.Cells.Item(65, 5).Value = "https://www.regionalstatistik.de/genesis//online?operation=table&code=33111-03-01-4&bypass=true&levelindex=1&levelid=1668501420853#abreadcrumb"
$ws.Cells.Item(65, 6).Value = ""

$ws.Cells.Item(66, 1).Value = "L_SP_1"
$ws.Cells.Item(66, 2).Value = "Q_VWGDL"
$ws.Cells.Item(66, 3).Value = "Investitionen, Anlagevermögen (eigene Berechnung auf Basis der Volkswirtschaftlichen Gesamtrechnungen der Länder)"
$ws.Cells.Item(66, 4).Value = "Investments, fixed assets (own calculation based on the Volkswirtschaftliche Gesamtrechnungen der Länder - only available in German)"
$ws.Cells.Item(66, 5).Value = "https://www.statistikportal.de/de/vgrdl/ergebnisse-laenderebene/investitionen-anlagevermoegen"
$ws.Cells.Item(66, 6).Value = "https://www.statistikportal.de/de/vgrdl/ergebnisse-laenderebene/investitionen-anlagevermoegen"

$ws.Cells.Item(67, 1).Value = "L_SP_2"
$ws.Cells.Item(67, 2).Value = "Q_VWGDL"
$ws.Cells.Item(67, 3).Value = "Bruttoinlandsprodukt, Bruttowertschöpfung (eigene Berechnung auf Basis der Volkswirtschaftlichen Gesamtrechnungen)"
$ws.Cells.Item(67, 4).Value = "Gross domestic product, gross value added (own calculation based on the Volkswirtschaftliche Gesamtrechnungen der Länder  - only available in German)"
$ws.Cells.Item(67, 5).Value = "https://www.statistikportal.de/de/vgrdl/ergebnisse-laenderebene/bruttoinlandsprodukt-bruttowertschoepfung"
$ws.Cells.Item(67, 6).Value = "https://www.statistikportal.de/de/vgrdl/ergebnisse-laenderebene/bruttoinlandsprodukt-bruttowertschoepfung"

$ws.Cells.Item(68, 1).Value = "L_SP_3"
$ws.Cells.Item(68, 2).Value = "Q_VWGDL"
$ws.Cells.Item(68, 3).Value = "Treibhausgasemissionen "
$ws.Cells.Item(68, 4).Value = "Greenhouse gas emissions (only available in German)"
$ws.Cells.Item(68, 5).Value = "https://www.statistikportal.de/de/ugrdl/ergebnisse/gase#alle-ergebnisse"
$ws.Cells.Item(68, 6).Value = "https://www.statistikportal.de/de/ugrdl/ergebnisse/gase#alle-ergebnisse"

$ws.Cells.Item(69, 1).Value = "L_SP_4"
$ws.Cells.Item(69, 2).Value = "Q_VWGDL"
$ws.Cells.Item(69, 3).Value = "Erwerbstätigenquoten"
$ws.Cells.Item(69, 4).Value = "Employment rates (only available in German)"
$ws.Cells.Item(69, 5).Value = "https://www.statistikportal.de/de/nachhaltigkeit/ergebnisse/ziel-8-menschenwuerdige-arbeit-und-wirtschaftswachstum"
$ws.Cells.Item(69, 6).Value = ""

$ws.Cells.Item(70, 1).Value = "L_SVWS"
$ws.Cells.Item(70, 2).Value = "Q_SVWS"
$ws.Cells.Item(70, 3).Value = "Forschung und Entwicklung in der Wirtschaft"
$ws.Cells.Item(70, 4).Value = "Research and development in the economy (only available in German)"
$ws.Cells.Item(70, 5).Value = "https://www.stifterverband.org/fue-facts-2020"
$ws.Cells.Item(70, 6).Value = ""

$ws.Cells.Item(71, 1).Value = "L_TI_1"
$ws.Cells.Item(71, 2).Value = "Q_TA"
$ws.Cells.Item(71, 3).Value = "Corruption Perceptions Index (nicht auf Deutsch verfügbar)"
$ws.Cells.Item(71, 4).Value = "Corruption Perceptions Index"
$ws.Cells.Item(71, 5).Value = "https://www.transparency.org/en/cpi/2020/index/nzl"
$ws.Cells.Item(71, 6).Value = "https://www.transparency.org/en/cpi/2020/index/nzl"

$ws.Cells.Item(72, 1).Value = "L_UBA_1"
$ws.Cells.Item(72, 2).Value = "Q_UBA"
$ws.Cells.Item(72, 3).Value = "Index der Luftschadstoff-Emissionen"
$ws.Cells.Item(72, 4).Value = "Index of air pollutant emissions"
$ws.Cells.Item(72, 5).Value = "https://www.umweltbundesamt.de/bild/index-der-luftschadstoff-emissionen"
$ws.Cells.Item(72, 6).Value = "https://www.umweltbundesamt.de/en/image/index-of-air-pollutant-emissions"

$ws.Cells.Item(73, 1).Value = "L_UBA_2"
$ws.Cells.Item(73, 2).Value = "Q_UBALAWA"
$ws.Cells.Item(73, 3).Value = "Messstellen an Flüssen mit Überschreitung des Orientierungswertes für Gesamtphosphor"
$ws.Cells.Item(73, 4).Value = "Sampling sites which exceeded the requirement for good status for total phosphorus in rivers"
$ws.Cells.Item(73, 5).Value = "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-eutrophierung-von-fluessen-durch-phosphor"
$ws.Cells.Item(73, 6).Value = "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-river-eutrophication-phosphorus"

$ws.Cells.Item(74, 1).Value = "L_UBA_3"
$ws.Cells.Item(74, 2).Value = "Q_UBA"
$ws.Cells.Item(74, 3).Value = "Nitrat im Grundwasser"
$ws.Cells.Item(74, 4).Value = "Nitrate in groundwater"
$ws.Cells.Item(74, 5).Value = "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-nitrat-im-grundwasser"
$ws.Cells.Item(74, 6).Value = "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-nitrate-in-groundwater"

$ws.Cells.Item(75, 1).Value = "L_UBA_4"
$ws.Cells.Item(75, 2).Value = "Q_UBA"
$ws.Cells.Item(75, 3).Value = "Marktanteile von Produkten mit staatlichen Umweltzeichen, nach Umsätzen gewichtet"
$ws.Cells.Item(75, 4).Value = "Weighted market shares by sales of products with official eco-labels"
$ws.Cells.Item(75, 5).Value = "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-umweltfreundlicher-konsum"
$ws.Cells.Item(75, 6).Value = "https://www.umweltbundesamt.de/en/indicator-environmentally-friendly-consumption"

$ws.Cells.Item(76, 1).Value = "L_UBA_5"
$ws.Cells.Item(76, 2).Value = "Q_UBA"
$ws.Cells.Item(76, 3).Value = "Emission der von der UN-Klimarahmenkonvention abgedeckten Treibhausgase"
$ws.Cells.Item(76, 4).Value = "Emission of greenhouse gases covered by the UN Framework Convention on Climate"
$ws.Cells.Item(76, 5).Value = "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-emission-von-treibhausgasen"
$ws.Cells.Item(76, 6).Value = "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-greenhouse-gas-emissions"

$ws.Cells.Item(77, 1).Value = "L_UBA_6"
$ws.Cells.Item(77, 2).Value = "Q_UBA"
$ws.Cells.Item(77, 3).Value = "Eutrophierung von Nord- und Ostsee durch Stickstoff"
$ws.Cells.Item(77, 4).Value = "Eutrophication of the North Sea / Baltic Sea by nitrogen"
$ws.Cells.Item(77, 5).Value = "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-eutrophierung-der-meere"
$ws.Cells.Item(77, 6).Value = "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-eutrophication-of-the-north-sea-baltic-sea"

$ws.Cells.Item(78, 1).Value = "L_UBA_7"
$ws.Cells.Item(78, 2).Value = "Q_UBA"
$ws.Cells.Item(78, 3).Value = "Anteil der Fläche empfindlicher Land-Ökosysteme mit Überschreitung der Belastungsgrenzen für Eutrophierung"
$ws.Cells.Item(78, 4).Value = "Proportion of vulnerable ecosystems where critical loads for eutrophication are exceeded"
$ws.Cells.Item(78, 5).Value = "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-eutrophierung-durch-stickstoff"
$ws.Cells.Item(78, 6).Value = "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-nitrogen-eutrophication"

# Step 3: fix auto row-height triggered by embedded newlines, restoring default row height
$ws.Rows.Item(33).AutoFit()
